$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain plain text,
# then restore the original (default) cell style so formatting is unaffected.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '70.063.92'
Set-TextValue $ws.Range("E2") '  -0.48%  '
# Row 3
Set-TextValue $ws.Range("D3") '3.499.87'
Set-TextValue $ws.Range("E3") '  -0.85%  '
# Row 4
Set-TextValue $ws.Range("E4") '  +0.18%  '
# Row 5
Set-TextValue $ws.Range("D5") '603.95'
Set-TextValue $ws.Range("E5") '  -0.70%  '
# Row 6
Set-TextValue $ws.Range("D6") '172.52'
Set-TextValue $ws.Range("E6") '  -0.45%  '
# Row 7
Set-TextValue $ws.Range("D7") '0.607'
Set-TextValue $ws.Range("E7") '  -1.71%  '
# Row 8
Set-TextValue $ws.Range("D8") '3.492.99'
Set-TextValue $ws.Range("E8") '  -0.79%  '
# Row 9
Set-TextValue $ws.Range("E9") '  +0.00%  '
# Row 10
Set-TextValue $ws.Range("D10") '0.194'
Set-TextValue $ws.Range("E10") '  -3.45%  '
# Row 11
Set-TextValue $ws.Range("E11") '  +6.74%  '
# Row 12
Set-TextValue $ws.Range("E12") '  +0.15%  '
# Row 13
Set-TextValue $ws.Range("E13") '  -3.08%  '
# Row 14
Set-TextValue $ws.Range("D14") '0.0000275'
Set-TextValue $ws.Range("E14") '  -2.05%  '
# Row 15
Set-TextValue $ws.Range("D15") '4.069.66'
Set-TextValue $ws.Range("E15") '  -0.45%  '
# Row 16
Set-TextValue $ws.Range("B16") 'Polkadot'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D16") '8.35'
Set-TextValue $ws.Range("E16") '  -0.97%  '
# Row 17
Set-TextValue $ws.Range("B17") 'BitcoinCash'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D17") '611.97'
Set-TextValue $ws.Range("E17") '  -2.40%  '
# Row 18
Set-TextValue $ws.Range("D18") '3.503.98'
Set-TextValue $ws.Range("E18") '  -0.46%  '
# Row 19
Set-TextValue $ws.Range("D19") '70.095.23'
Set-TextValue $ws.Range("E19") '  -0.16%  '
# Row 20
Set-TextValue $ws.Range("E20") '  +0.89%  '
# Row 21
Set-TextValue $ws.Range("D21") '17.53'
Set-TextValue $ws.Range("E21") '  +0.82%  '
# Row 22
Set-TextValue $ws.Range("D22") '0.877'
Set-TextValue $ws.Range("E22") '  -1.19%  '
# Row 23
Set-TextValue $ws.Range("D23") '9.06'
Set-TextValue $ws.Range("E23") '  -9.29%  '
# Row 24
Set-TextValue $ws.Range("D24") '98.73'
Set-TextValue $ws.Range("E24") '  +2.21%  '
# Row 25
Set-TextValue $ws.Range("E25") '  -1.87%  '
# Row 26
Set-TextValue $ws.Range("E26") '  -3.96%  '
# Row 27
Set-TextValue $ws.Range("E27") '  -0.14%  '
# Row 28
Set-TextValue $ws.Range("D28") '2.55'
Set-TextValue $ws.Range("E28") '  -2.14%  '
# Row 29
Set-TextValue $ws.Range("D29") '33.97'
Set-TextValue $ws.Range("E29") '  +1.97%  '
# Row 30
Set-TextValue $ws.Range("D30") '8.99'
Set-TextValue $ws.Range("E30") '  -2.82%  '
# Row 31
Set-TextValue $ws.Range("D31") '8.02'
Set-TextValue $ws.Range("E31") '  -5.25%  '
# Row 32
Set-TextValue $ws.Range("E32") '  -4.75%  '
# Row 33
Set-TextValue $ws.Range("B33") 'Mantle'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D33") '1.28'
Set-TextValue $ws.Range("E33") '  -4.87%  '
# Row 34
Set-TextValue $ws.Range("B34") 'Bittensor'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D34") '628.75'
Set-TextValue $ws.Range("E34") '  +12.06%  '
# Row 35
Set-TextValue $ws.Range("E35") '  -3.42%  '
# Row 36
Set-TextValue $ws.Range("D36") '0.0993'
Set-TextValue $ws.Range("E36") '  -2.59%  '
# Row 37
Set-TextValue $ws.Range("E37") '  -0.83%  '
# Row 38
Set-TextValue $ws.Range("D38") '3.49'
Set-TextValue $ws.Range("E38") '  -1.88%  '
# Row 39
Set-TextValue $ws.Range("E39") '  +5.64%  '
# Row 40
Set-TextValue $ws.Range("D40") '56.71'
Set-TextValue $ws.Range("E40") '  -0.90%  '
# Row 41
Set-TextValue $ws.Range("E41") '  +0.42%  '
# Row 42
Set-TextValue $ws.Range("D42") '0.143'
Set-TextValue $ws.Range("E42") '  +0.39%  '
# Row 43
Set-TextValue $ws.Range("D43") '3.357.71'
Set-TextValue $ws.Range("E43") '  +0.43%  '
# Row 44
Set-TextValue $ws.Range("D44") '0.0₃0730'
Set-TextValue $ws.Range("E44") '  +2.02%  '
# Row 45
Set-TextValue $ws.Range("E45") '  -5.99%  '
# Row 46
Set-TextValue $ws.Range("E46") '  -3.77%  '
# Row 47
Set-TextValue $ws.Range("D47") '31.84'
Set-TextValue $ws.Range("E47") '  -3.98%  '
# Row 48
Set-TextValue $ws.Range("D48") '2.54'
Set-TextValue $ws.Range("E48") '  -3.87%  '
# Row 49
Set-TextValue $ws.Range("E49") '  +0.50%  '
# Row 50
Set-TextValue $ws.Range("D50") '132.90'
Set-TextValue $ws.Range("E50") '  -1.22%  '
